# Generate Report for Archive
#
# The localization status for this report's single row moves from
# "Ready for handoff" to "In Translation". That text lives in one shared
# string used by the Status column on every sheet: Overview!E2:F2
# (zh-cn / de-de columns) and C2 on each of the per-language sheets.
# Narrow the now-shorter Status column(s) to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = "In Translation"
$overview.Range("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C:C").ColumnWidth = 12.5
